$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column G
$ws.Range("G1").Value = "s1cDNAProtocol"

# Fill G2:G19 with the protocol value
$ws.Range("G2:G19").Value = "E7420L"
$ws.Range("G2:G19").WrapText = $true

# Update selection to match the target state
$null = $ws.Range("G2:G19").Select()
